$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.249.08'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").Value = '3.429.70'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.483'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.414'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '4.017.14'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.05%  '
$ws.Range("D15").Value = '3.401.91'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '62.221.88'
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.568'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '3.556.22'
$ws.Range("E25").Value = '  -1.73%  '
$ws.Range("E26").Value = '  -3.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '167.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").Value = '3.462.30'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0780'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.778'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.24%  '
$ws.Range("D47").Value = '2.544.91'
$ws.Range("E47").Value = '  -3.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.17%  '
